$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for all data rows (2 through 189)
# from serial date 45205 (2023-10-06) to 45206 (2023-10-07).
$ws.Range("C2:C189").Value = 45206
